# Update "location cleanup" sheet: strip the leading "case_when" token from
# the formula-fragment text stored in column C (shared across C3:C174), so
# that the generated R snippets in column F read "(chicago_crime$location_description == '..."
# instead of "case_when(chicago_crime$location_description == '...".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("location cleanup")

$ws.Range("C3:C174").Value = "(chicago_crime`$location_description == '"
